# Apply crypto price/volume updates for Sun Apr  2 03:46:02 UTC 2023 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting for numeric-looking price values (column D is stored as text)
$dCells = @("D2", "D3", "D5", "D6", "D7", "D8", "D9", "D10", "D12", "D13", "D15", "D16", "D17", "D18", "D20", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D31", "D32", "D33", "D34", "D36", "D37", "D38", "D39", "D40", "D41", "D45", "D46", "D47", "D48", "D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "28.539.84"
$ws.Range("D3").Value = "1.820.33"
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "315.90"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").Value = "0.5153"
$ws.Range("E7").Value = "  -3.53%  "
$ws.Range("D8").Value = "0.3876"
$ws.Range("E8").Value = "  -2.68%  "
$ws.Range("D9").Value = "0.08404"
$ws.Range("E9").Value = "  +7.57%  "
$ws.Range("D10").Value = "41.89"
$ws.Range("E10").Value = "  -0.51%  "
$ws.Range("E11").Value = "  -0.61%  "
$ws.Range("D12").Value = "6.404"
$ws.Range("D13").Value = "21.05"
$ws.Range("E13").Value = "  -0.18%  "
$ws.Range("E14").Value = "  -0.05%  "
$ws.Range("D15").Value = "7.510"
$ws.Range("E15").Value = "  -0.83%  "
$ws.Range("D16").Value = "1.813.46"
$ws.Range("E16").Value = "  -0.91%  "
$ws.Range("D17").Value = "0.00001135"
$ws.Range("E17").Value = "  +3.86%  "
$ws.Range("D18").Value = "92.84"
$ws.Range("E18").Value = "  -0.49%  "
$ws.Range("E19").Value = "  +2.01%  "
$ws.Range("D20").Value = "17.78"
$ws.Range("E20").Value = "  -0.21%  "
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").Value = "6.093"
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("D23").Value = "28.579.05"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").Value = "11.44"
$ws.Range("E24").Value = "  +2.03%  "
$ws.Range("D25").Value = "2.275"
$ws.Range("E25").Value = "  +1.51%  "
$ws.Range("D26").Value = "21.10"
$ws.Range("E26").Value = "  +1.16%  "
$ws.Range("D27").Value = "159.02"
$ws.Range("E27").Value = "  +1.18%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "2.422"
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("B29").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C29").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D29").Value = "2.023.58"
$ws.Range("E29").Value = "  -0.82%  "
$ws.Range("E30").Value = "  +0.49%  "
$ws.Range("D31").Value = "1.096"
$ws.Range("E31").Value = "  -4.84%  "
$ws.Range("D32").Value = "0.1080"
$ws.Range("E32").Value = "  -4.02%  "
$ws.Range("D33").Value = "5.757"
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("D34").Value = "0.07536"
$ws.Range("E34").Value = "  +3.00%  "
$ws.Range("E35").Value = "  +0.81%  "
$ws.Range("D36").Value = "0.2226"
$ws.Range("E36").Value = "  -1.73%  "
$ws.Range("D37").Value = "0.02365"
$ws.Range("E37").Value = "  +0.61%  "
$ws.Range("D38").Value = "5.208"
$ws.Range("E38").Value = "  -0.11%  "
$ws.Range("D39").Value = "8.733"
$ws.Range("E39").Value = "  -3.00%  "
$ws.Range("D40").Value = "0.6338"
$ws.Range("E40").Value = "  +0.58%  "
$ws.Range("D41").Value = "11.25"
$ws.Range("E41").Value = "  -1.31%  "
$ws.Range("E42").Value = "  -0.75%  "
$ws.Range("E43").Value = "  +0.46%  "
$ws.Range("E44").Value = "  +0.32%  "
$ws.Range("D45").Value = "3.770"
$ws.Range("E45").Value = "  +1.48%  "
$ws.Range("D46").Value = "0.5926"
$ws.Range("E46").Value = "  -0.12%  "
$ws.Range("D47").Value = "125.72"
$ws.Range("E47").Value = "  +0.14%  "
$ws.Range("D48").Value = "1.991"
$ws.Range("E48").Value = "  -0.33%  "
$ws.Range("E49").Value = "  +0.36%  "
$ws.Range("E50").Value = "  +0.53%  "
$ws.Range("D51").Value = "74.39"
$ws.Range("E51").Value = "  -0.35%  "
